{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The edit:\n//   1. Right after the page's title paragraph (\"Play Captain's Quest\n//      Treasure Island Free - Review 2021\", Heading 1), insert a new\n//      paragraph with a bold \"Meta description\" label followed by the\n//      SEO summary text.\n//   2. At the bottom of the document, remove the paragraph that had\n//      duplicated the bold title text, and replace the text of the\n//      trailing italic paragraph (previously the SEO summary) with a\n//      new \"feature image\" prompt, keeping its italic formatting.\n\nconst body = context.document.body;\n\nconst titleText = \"Play Captain's Quest Treasure Island Free - Review 2021\";\nconst metaRestText =\n  \": Get the ultimate free play review of Captain's Quest Treasure Island Slot, a 5-reel game with 10 paylines and a high volatility of 96% RTP.\";\nconst featureImageText =\n  \"For the feature image, create a cartoon-style design featuring a Maya warrior with glasses who is looking happy and satisfied. The design should include elements of the game, such as a ship sailing the Caribbean Sea, a deserted island where the treasure is hidden, and symbols of the game like the poker card suits, the helm, and the treasure. The background of the image should be blue with a pirate-themed border, and the game's name \\\"Captain's Quest Treasure Island\\\" should be prominently displayed. Make sure the image is bright and eye-catching, with lots of detail to entice players to try out the game.\";\n\n// ---------------------------------------------------------------------\n// 1) Insert the \"Meta description: ...\" paragraph right after the\n//    title paragraph (the Heading 1 at the very top of the document).\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\n\nconst metaPara = titlePara.insertParagraph(\"\", Word.InsertLocation.after);\nmetaPara.style = \"Normal\";\nawait context.sync();\n\nmetaPara.insertText(\"Meta description\" + metaRestText, Word.InsertLocation.end);\nawait context.sync();\n\n// Bold just the \"Meta description\" label, leaving the rest of the\n// sentence in regular weight.\nconst metaLabelHits = metaPara.search(\"Meta description\", { matchCase: true });\nmetaLabelHits.load(\"items\");\nawait context.sync();\nmetaLabelHits.items[0].font.bold = true;\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) At the bottom of the document: delete the paragraph duplicating\n//    the bold title text, and replace the italic paragraph's text\n//    with the new feature-image prompt (formatting stays italic).\n// ---------------------------------------------------------------------\nconst titleHits = body.search(titleText, { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\n\n// The last match is the duplicate bold paragraph near the end of the\n// document (the first match is the Heading 1 title itself).\nconst duplicateTitleRange = titleHits.items[titleHits.items.length - 1];\nconst duplicateTitleParas = duplicateTitleRange.paragraphs;\nduplicateTitleParas.load(\"items\");\nawait context.sync();\nduplicateTitleParas.items[0].delete();\nawait context.sync();\n\n// The last paragraph in the document is the italic SEO-summary one;\n// replace its text while keeping the paragraph (and its italic run).\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\nconst lastPara = allParagraphs.items[allParagraphs.items.length - 1];\nlastPara.insertText(featureImageText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# The edit:\n#   1. Right after the page's title paragraph (\"Play Captain's Quest\n#      Treasure Island Free - Review 2021\", Heading 1), insert a new\n#      paragraph with a bold \"Meta description\" label followed by the\n#      SEO summary text.\n#   2. At the bottom of the document, remove the paragraph that had\n#      duplicated the bold title text, and replace the text of the\n#      trailing italic paragraph (previously the SEO summary) with a\n#      new \"feature image\" prompt, keeping its italic formatting.\n\n$d = $word.ActiveDocument\n\n$titleText = \"Play Captain's Quest Treasure Island Free - Review 2021\"\n$metaFullText = \"Meta description: Get the ultimate free play review of Captain's Quest Treasure Island Slot, a 5-reel game with 10 paylines and a high volatility of 96% RTP.\"\n$metaLabelLength = 16  # length of \"Meta description\"\n$featureImageText = \"For the feature image, create a cartoon-style design featuring a Maya warrior with glasses who is looking happy and satisfied. The design should include elements of the game, such as a ship sailing the Caribbean Sea, a deserted island where the treasure is hidden, and symbols of the game like the poker card suits, the helm, and the treasure. The background of the image should be blue with a pirate-themed border, and the game's name `\"Captain's Quest Treasure Island`\" should be prominently displayed. Make sure the image is bright and eye-catching, with lots of detail to entice players to try out the game.\"\n\n# ---------------------------------------------------------------------\n# 1) Insert the \"Meta description: ...\" paragraph right after the\n#    title paragraph (the first paragraph / Heading 1 of the document).\n# ---------------------------------------------------------------------\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaPara.Style = \"Normal\"\n\n$metaRange = $metaPara.Range\n$metaRange.Text = $metaFullText\n\n# Bold just the \"Meta description\" label, leaving the rest plain.\n$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $metaLabelLength)\n$boldRange.Bold = 1\n\n# ---------------------------------------------------------------------\n# 2) At the bottom of the document: delete the paragraph duplicating\n#    the bold title text, and replace the italic paragraph's text\n#    with the new feature-image prompt (formatting stays italic).\n# ---------------------------------------------------------------------\n$searchRange = $d.Range($metaRange.End, $d.Content.End)\n$searchRange.Find.MatchCase = $true\n$searchRange.Find.Text = $titleText\n$searchRange.Find.Execute() | Out-Null\n\n$duplicateTitlePara = $searchRange.Paragraphs(1)\n$duplicateTitlePara.Range.Delete()\n\n$count = $d.Paragraphs.Count\n$italicPara = $d.Paragraphs($count)\n$italicRange = $italicPara.Range\n$italicRange.MoveEnd(1, -1)\n$italicRange.Text = $featureImageText\n"}
